$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 471
$hyperlinkCols = 19,20,21,22,23,24,25

for ($row = 2; $row -le $lastRow; $row++) {
    $ws.Cells.Item($row, 3).Value = 45186

    if ($row -le 23) {
        $designation = $ws.Cells.Item($row, 1).Value()
        foreach ($col in $hyperlinkCols) {
            $cell = $ws.Cells.Item($row, $col)
            if ($cell.HasFormula) {
                $formula = $cell.Formula
                if ($formula.EndsWith(')')) {
                    $newFormula = $formula.Substring(0, $formula.Length - 1) + ', "' + $designation + '")'
                    $cell.Formula = $newFormula
                }
            }
        }
    }
}
